$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F5").Value = 5251
$ws.Range("F6").Value = 5251
$ws.Range("F12").Value = 743
$ws.Range("F13").Value = 5186
$ws.Range("F15").Value = 74
$ws.Range("F17").Value = 926
$ws.Range("F18").Value = 927
$ws.Range("F22").Value = 3906
$ws.Range("F24").Value = 3829
$ws.Range("F30").Value = 210
$ws.Range("F37").Value = 6808
$ws.Range("F42").Value = 61
$ws.Range("F43").Value = 1386
$ws.Range("F47").Value = 2316
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F25").Value = 816
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F7").Value = 5251
$ws.Range("F8").Value = 5251
$ws.Range("F14").Value = 743
$ws.Range("F15").Value = 5186
$ws.Range("F17").Value = 74
$ws.Range("F19").Value = 930
$ws.Range("F20").Value = 930
$ws.Range("F24").Value = 3906
$ws.Range("F25").Value = 3829
$ws.Range("F30").Value = 210
$ws.Range("F37").Value = 6808
$ws.Range("F43").Value = 61
$ws.Range("F44").Value = 1386
$ws.Range("F47").Value = 2316
